# 20180103_LN3 Deep Completed (House14)_newcleaned.xlsx
# - Renamed worksheet "LN2" -> "LN3" (the sheet's _FilterDatabase defined
#   name is range-scoped to the sheet itself and updates automatically).
# - Cleared the AutoFilter criteria that had been restricting column H
#   ("Benthic Cat.") to "HC" only, which re-reveals every previously
#   filtered-out row without touching the autofilter range itself.
# - Moved the active selection down to C314 (near the bottom of the data)
#   to reflect where the author was working after un-filtering.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LN2")

# Rename the sheet; Excel keeps the sheet-scoped _xlnm._FilterDatabase
# defined name (and any formulas referencing the sheet) in sync.
$ws.Name = "LN3"

# Remove the "HC"-only filter on the "Benthic Cat." column (column 8)
# while keeping the AutoFilter range A1:T311 itself intact.
$ws.ShowAllData()

# Move to where the editing continued, near the end of the filtered data.
$ws.Activate()
$ws.Range("C314").Select()
